$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "compulsarySchooling" program as row 48 -------------------
# New shared-strings entries are appended to sharedStrings.xml in the order
# they are first written, so the string-valued cells are populated in the
# J, I, F, A, B order to reproduce the target shared-string indices.
$ws.Range("J48").Value = "https://ideas.repec.org/a/bla/scandj/v116y2014i3p878-907.html"
$ws.Range("I48").Value = "Piopiunik (2014)"
$ws.Range("F48").Value = "Between 1946 and 1969 all German federal states extended the length of the least academic school track `"Hauptschule`" in Germany by one year. Piopiunik (2014) finds that this reform improved educational outcomes of sons whose mothers were affected by the compulsary schooling reform."
$ws.Range("A48").Value = "compulsarySchooling"
$ws.Range("B48").Value = "Compulsary Schooling"
$ws.Range("C48").Value = 1967
$ws.Range("D48").Value = "Education"
$ws.Range("E48").Value = 16

# Row 48 needs extra height to show the wrapped description text.
$ws.Rows.Item(48).RowHeight = 105

# Turn the URL literal in J48 into a clickable hyperlink (adds a new
# relationship + <hyperlink> entry, same as the other rows in column J).
$ws.Hyperlinks.Add($ws.Range("J48"), "https://ideas.repec.org/a/bla/scandj/v116y2014i3p878-907.html") | Out-Null

# Hyperlinks.Add() stamps the cell with a brand new "Link" style; reapply
# the formatting already used by the other hyperlinked cells in column J
# (e.g. J47) so J48 re-uses that existing cell style instead of a new one.
$ws.Range("J47").Copy()
$ws.Range("J48").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Update the view/selection shown when the sheet is reopened ------------
$ws.Range("L48").Select() | Out-Null
